$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (Table2) currently spans A1:E80 (79 data rows). Add one more row
# via the ListObject so the table range / autofilter grow automatically.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()
$rng = $newRow.Range

# Question / Difficulty / Pattern / Notes
$rng.Item(1, 1).Value = "73. Set Matrix Zeroes"
$rng.Item(1, 3).Value = "Math"
$rng.Item(1, 4).Value = "In the brute force, start with a copy matrix. Read the input, but update the copy, as we do not want it to dynamically update. For an improvement, replace the copy matrix with a row array and column array to mark which rows and columns get set to zero. The optimal, for O(1) space, we put the row array and column array inside the input matrix itself to perform in place, but need 1 extra dedicated variable for the first cell where rows and columns overlap. The reason we can overlap, is because we read the input before we overwrite, so the computation is accurate."

# Difficulty "Medium" should reuse the same yellow-fill style as the other
# Medium rows (e.g. B80) instead of getting a brand new cell style, so copy
# the format from the row above first, then set the value.
$ws.Range("B80").Copy()
$rng.Item(1, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$rng.Item(1, 2).Value = "Medium"

# Link column: add the real hyperlink, then restore the shared "Hyperlink"
# cell style so it matches the rest of column E.
$linkCell = $rng.Item(1, 5)
$ws.Hyperlinks.Add($linkCell, "https://leetcode.com/problems/set-matrix-zeroes/solutions/2525398/all-approaches-from-brute-force-to-optimal-with-easy-explanation/ ", "", "", "https://leetcode.com/problems/set-matrix-zeroes/solutions/2525398/all-approaches-from-brute-force-to-optimal-with-easy-explanation/ ")
$linkCell.Style = "Hyperlink"

# Match the author's final cursor position / scrolled view.
$ws.Range("D87").Select() | Out-Null

Write-Host $ws.UsedRange.Address()
